# Update odds values on Sheet1 to reflect the latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("G8").Value = 1.87
$ws.Range("N8").Value = 1.95
$ws.Range("O8").Value = 1.85

# Row 9
$ws.Range("G9").Value = 4.5
$ws.Range("H9").Value = 4
$ws.Range("L9").Value = 1.25
$ws.Range("M9").Value = 3.75
$ws.Range("P9").Value = 1.3
$ws.Range("T9").Value = 13
$ws.Range("U9").Value = 23
$ws.Range("X9").Value = 34
$ws.Range("Z9").Value = 13
$ws.Range("AE9").Value = 8

# Row 10
$ws.Range("G10").Value = 2.35
$ws.Range("J10").Value = 1.06
$ws.Range("K10").Value = 10
$ws.Range("L10").Value = 1.3
$ws.Range("M10").Value = 3.4
$ws.Range("N10").Value = 2
$ws.Range("O10").Value = 1.8
$ws.Range("P10").Value = 1.37
$ws.Range("Q10").Value = 2.75
$ws.Range("R10").Value = 1.73
$ws.Range("S10").Value = 2
$ws.Range("T10").Value = 8
$ws.Range("U10").Value = 11
$ws.Range("X10").Value = 19
$ws.Range("Y10").Value = 29
$ws.Range("Z10").Value = 10
$ws.Range("AB10").Value = 15
$ws.Range("AC10").Value = 51
$ws.Range("AD10").Value = 251
$ws.Range("AE10").Value = 9.5
$ws.Range("AJ10").Value = 34

# Row 11
$ws.Range("L11").Value = 1.29
$ws.Range("M11").Value = 3.5
$ws.Range("N11").Value = 1.95
$ws.Range("O11").Value = 1.85

# Row 13
$ws.Range("J13").Value = 1.04
$ws.Range("L13").Value = 1.25
$ws.Range("N13").Value = 1.8
$ws.Range("O13").Value = 2

# Row 14
$ws.Range("J14").Value = 1.03
$ws.Range("L14").Value = 1.19
$ws.Range("M14").Value = 4
